$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely: every subsequent row (old row 3..82) shifts up by one,
# so the old "2025-10-13" row disappears and the trailing "2026-01-01" row
# (old row 82) is gone too since the used range shrinks by one row.
$ws.Rows.Item(2).Delete()

# The two newly-shifted-up rows (now rows 2 and 3, dated 2025-10-14 and
# 2025-10-15) lose their "No video indexed" / "Video indexed" counts -
# blank them out.
$ws.Range("B2:C3").Value = ""
